$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "ingresso disponível" use-case entries to "Tipo de Ingresso"
# wording (column B, rows 27-30). Edit order matches the shared-string
# append order seen in the target workbook (Exibir, Excluir, Editar,
# Cadastrar).
$ws.Range("B28").Value = "Exibir Tipo de Ingresso"
$ws.Range("B30").Value = "Excluir Tipo de Ingresso"
$ws.Range("B29").Value = "Editar Tipo de Ingresso"
$ws.Range("B27").Value = "Cadastrar Tipo de ingresso"

# Row 17 picked up an explicit (custom) height during the edit session.
$ws.Rows.Item(17).RowHeight = 31.5

# Reflect the cursor ending up on A35 after the edits.
$ws.Range("A35").Select()
